$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after the current "Pseudotime" row (row 35) so the
# single summary row can be split into three (Pseudotime_1/_2/_3). This
# pushes everything from the old row 36 onward down by two rows.
$ws.Rows.Item(36).EntireRow.Insert()
$ws.Rows.Item(36).EntireRow.Insert()

# Row 35: Pseudotime -> Pseudotime_1, with new statistics
$ws.Range("A35").Value = "Pseudotime_1"
$ws.Range("B35").Value = "  5.98 (3.85)  "
$ws.Range("C35").Value = " 4.51 (3.25)  "
$ws.Range("D35").Value = " <0.001  "

# Row 36 (new): Pseudotime_2
$ws.Range("A36").Value = "Pseudotime_2"
$ws.Range("B36").Value = "  4.67 (2.74)  "
$ws.Range("C36").Value = " 4.54 (3.03)  "
# This p-value text ("  0.669  ") looks like a plain number, so force the
# cell to Text format first or Excel will silently coerce it to a Double.
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "  0.669  "

# Row 37 (new): Pseudotime_3
$ws.Range("A37").Value = "Pseudotime_3"
$ws.Range("B37").Value = "  4.40 (3.03)  "
$ws.Range("C37").Value = " 4.08 (2.25)  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "  0.254  "
